$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" sheet entirely
$wsDel = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDel.Delete() | Out-Null

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"
